{"js": "// Update the attendance table (page 3 / \"crear el curso\") with the new\n// student data for the two data rows of the single table in the body.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// Row 1 (first student): name (with line break), document number, program.\nconst row1Name = table.getCell(1, 1);\nrow1Name.value = \"Samuel Andres Celis\\u000bLizcano\";\n\nconst row1Doc = table.getCell(1, 2);\nrow1Doc.value = \"1091964042\";\n\nconst row1Program = table.getCell(1, 3);\nrow1Program.value = \"ING de Sistemas\";\n\n// Row 2 (second student): name, document number, program, attendance.\nconst row2Name = table.getCell(2, 1);\nrow2Name.value = \"pablo Bb\";\n\nconst row2Doc = table.getCell(2, 2);\nrow2Doc.value = \"1004922828\";\n\nconst row2Program = table.getCell(2, 3);\nrow2Program.value = \"ING de Sistemas\";\n\nconst row2Attendance = table.getCell(2, 4);\nrow2Attendance.value = \"NO ASISTI\u00d3\";\n\nawait context.sync();\n", "ps1": "# Update the attendance table (page 3 / \"crear el curso\") with the new\n# student data for the two data rows of the single table in the body.\n$d = $word.ActiveDocument\n$t = $d.Tables(1)\n\n# Row 2 in COM 1-based indexing = first student data row.\n$t.Cell(2, 2).Range.Text = \"Samuel Andres Celis\" + [char]11 + \"Lizcano\"\n$t.Cell(2, 3).Range.Text = \"1091964042\"\n$t.Cell(2, 4).Range.Text = \"ING de Sistemas\"\n\n# Row 3 in COM 1-based indexing = second student data row.\n$t.Cell(3, 2).Range.Text = \"pablo Bb\"\n$t.Cell(3, 3).Range.Text = \"1004922828\"\n$t.Cell(3, 4).Range.Text = \"ING de Sistemas\"\n$t.Cell(3, 5).Range.Text = \"NO ASISTI\u00d3\"\n"}
